# Edit script: add evaluator_partial_correctness column, update o_10 data, add o_20 and o_20_jumbled sheets
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add new header column E on sheet1 (o_10), copying style from D1 ---
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"

# --- Update existing data row on sheet1 (o_10) ---
$p10 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
is the following a valid eulerian graph, if traversal is started from A?
   A B C D E F G H I J
 A 0 0 0 0 0 0 0 0 1 0
 B 0 0 1 0 0 1 1 0 1 0
 C 0 1 0 1 0 0 0 0 0 1
 D 0 0 1 0 0 0 0 0 1 0
 E 0 0 0 0 0 1 1 0 0 1
 F 0 1 0 0 1 0 0 0 1 0
 G 0 1 0 0 1 0 0 1 1 0
 H 0 0 0 0 0 0 1 0 0 0
 I 1 1 0 1 0 1 1 0 0 0
 J 0 0 1 0 1 0 0 0 0 0
    
"@
$ws1.Range("A2").Value = $p10
$ws1.Range("B2").Value = "This is not a valid eulerian graph"
$r10 = @"
To determine if the given graph is a valid Eulerian graph when starting traversal from node A, we need to check if every node in the graph has an even degree. 
The degree of a node is the number of edges that connect to it. In the given adjacency matrix, the value in the Mth row and Nth column represents the connection between nodes M and N.
Let's calculate the degree of each node:
Node A: 1 (There is a connection with node I)
Node B: 5 (Connections with nodes C, F, G, H, and I)
Node C: 4 (Connections with nodes B, D, J)
Node D: 2 (Connections with nodes C, I)
Node E: 3 (Connections with nodes F, G, J)
Node F: 3 (Connections with nodes B, E, I)
Node G: 4 (Connections with nodes B, E, F, H)
Node H: 1 (Connection with node G)
Node I: 5 (Connections with nodes B, D, E, F, G)
Node J: 3 (Connections with nodes C, E, I)
We can see that nodes A, C, D, H, and J have an odd degree while all other nodes have an even degree. In Eulerian graphs, all nodes except for 2 can have an odd degree. 
Therefore, the given graph is not a valid Eulerian graph when starting traversal from node A.
"@
$ws1.Range("C2").Value = $r10
$ws1.Range("D2").Value = "Wrong"
$ws1.Range("E2").Value = "N/A"
$ws1.Rows.Item(2).AutoFit()

# --- Create sheet o_20 after o_10, copying header formatting ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"
$ws1.Range("A1:E1").Copy($ws2.Range("A1"))

$p20 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
is the following a valid eulerian graph, if traversal is started from 6?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0
 B 1 0 1 0 0 1 1 1 0 0 0 1 0 0 0 0 0 0 0 0
 C 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 1 1 0 1 0 1 0 0 0 1 1 0 0 0 0 0 0 0
 F 0 1 0 1 1 0 1 0 1 0 0 0 0 0 0 0 0 0 1 0
 G 0 1 0 0 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 H 0 1 0 1 1 0 1 0 1 0 1 0 0 1 0 0 0 1 0 0
 I 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 1 0 0 1 1
 J 1 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 1 0
 K 0 0 0 0 0 0 1 1 0 1 0 1 0 0 0 0 0 0 0 0
 L 0 1 0 0 1 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 1 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0
 N 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0
 O 1 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 1 0 0 0
 P 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0
 R 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0
 S 0 0 0 0 0 1 0 0 1 1 0 0 0 0 0 0 0 0 0 1
 T 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0
    
"@
$ws2.Range("A2").Value = $p20
$ws2.Range("B2").Value = "This is a valid eulerian graph"
$r20 = @"
To determine if this graph is a valid Eulerian graph, we need to check if every vertex has an even degree. 
We can start by creating a list to store the degree of each vertex:
A: degree = 1
B: degree = 8
C: degree = 3
D: degree = 3
E: degree = 4
F: degree = 5
G: degree = 3
H: degree = 6
I: degree = 5
J: degree = 4
K: degree = 4
L: degree = 4
M: degree = 3
N: degree = 2
O: degree = 4
P: degree = 2
Q: degree = 2
R: degree = 2
S: degree = 4
T: degree = 3
As we can see, vertex A has an odd degree, which means this graph is not a valid Eulerian graph. Therefore, it is not possible to start traversal from vertex 6.
"@
$ws2.Range("C2").Value = $r20
$ws2.Range("D2").Value = "Wrong"
$ws2.Range("E2").Value = "N/A"
$ws2.Rows.Item(2).AutoFit()

# --- Create sheet o_20_jumbled after o_20, copying header formatting ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"
$ws1.Range("A1:E1").Copy($ws3.Range("A1"))

$p20j = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
is the following a valid eulerian graph, if traversal is started from A?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 C 1 0 0 0 0 0 0 1 0 1 0 1 0 0 0 0 1 0 0 0
 D 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 1 1 0
 F 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 I 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1
 J 0 0 1 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 1 0 0 0 0 0 0 1 1 0 0 1 0 1 0 0 0 0
 M 0 0 0 0 0 1 1 0 0 0 0 0 0 0 0 1 0 0 0 0
 N 1 1 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0
 O 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 P 0 0 0 0 1 0 0 0 1 0 0 1 1 0 0 0 0 0 0 0
 Q 0 0 1 0 0 0 0 0 0 0 0 0 0 1 1 0 0 1 0 1
 R 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 S 0 0 0 1 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 1
 T 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1 0 1 0
    
"@
$ws3.Range("A2").Value = $p20j
$ws3.Range("B2").Value = "This is not a valid eulerian graph"
$r20j = @"
To determine if the given graph is valid eulerian graph, we need to check if the graph is connected and if every vertex has an even degree.
First, let's check for connectivity. We can do this by performing a Depth-First Search (DFS) or Breadth-First Search (BFS) starting from vertex A and checking if we can reach every other vertex.
Performing a BFS from vertex A, we can reach all other vertices. Therefore, the graph is connected.
Next, we need to check if every vertex has an even degree. We can do this by counting the number of ones in each row of the adjacency matrix and checking if the count is even for every vertex.
Counting the ones for each row:
A: 2
B: 1
C: 8
D: 1
E: 5
F: 1
G: 1
H: 2
I: 4
J: 3
K: 1
L: 6
M: 3
N: 6
O: 2
P: 4
Q: 6
R: 2
S: 4
T: 3
Since not every vertex has an even degree, the graph is not a valid eulerian graph.
Therefore, the given graph is not a valid eulerian graph if traversal is started from vertex A.
"@
$ws3.Range("C2").Value = $r20j
$ws3.Range("D2").Value = "Wrong"
$ws3.Range("E2").Value = "N/A"
$ws3.Rows.Item(2).AutoFit()

$ws1.Select()
Write-Output "done"
